# Apply "Fixed the bullet box friction test." edit.
# Mark several tests as Fixed and clear their associated error-message cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2: Test_BoxFriction (Bullet, template mismatch #1) -> now Fixed
$ws.Range("B2").Value = "Fixed"
$ws.Range("E2").ClearContents()

# Row 3: Test_BoxFriction (Bullet, template mismatch #2) -> now Fixed
$ws.Range("B3").Value = "Fixed"
$ws.Range("E3").ClearContents()

# Row 5: Test_ContactSensors -> now Fixed
$ws.Range("B5").Value = "Fixed"
$ws.Range("D5").ClearContents()

# Row 6: Test_DeleteNodesBeforeOpenChart (already Fixed) -> clear stale error note
$ws.Range("D6").ClearContents()

# Row 7: Test_Distance -> now Fixed
$ws.Range("B7").Value = "Fixed"
$ws.Range("D7").ClearContents()

# Row 25: Test_Torus -> now Fixed
$ws.Range("B25").Value = "Fixed"
$ws.Range("D25").ClearContents()

# Update the saved view: scroll back to top-left and select E4
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("E4").Select()
